# Prevent resizing, center window
# (Functionally: a new top leaderboard entry "Hello" with score 4 is
#  inserted at row 34, and the existing entries that were in rows 34-40
#  shift down into rows 35-41. The sheet does not grow; the data that
#  was in row 41 is dropped off the bottom of the table.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current values of rows 34-40 (before they get overwritten)
# so they can be written back out shifted down into rows 35-41.
$sourceRows = 34..40
$savedValues = @{}
foreach ($r in $sourceRows) {
    $savedValues[$r] = @(
        $ws.Range("A$r").Value2,
        $ws.Range("B$r").Value2,
        $ws.Range("C$r").Value2,
        $ws.Range("D$r").Value2,
        $ws.Range("E$r").Value2
    )
}

# Write the shifted-down rows (old row N -> new row N+1), starting from
# the bottom so we never overwrite a source row before it's been read.
foreach ($r in ($sourceRows | Sort-Object -Descending)) {
    $vals = $savedValues[$r]
    $target = $r + 1
    $ws.Range("A$target").Value = $vals[0]
    $ws.Range("B$target").Value = $vals[1]
    $ws.Range("C$target").Value = $vals[2]
    $ws.Range("D$target").Value = $vals[3]
    $ws.Range("E$target").Value = $vals[4]
}

# Insert the new leaderboard entry into row 34.
$ws.Range("A34").Value = "Teacher"
$ws.Range("B34").Value = "Hello"
$ws.Range("C34").Value = "N/A"
$ws.Range("D34").Value = "N/A"
$ws.Range("E34").Value = 4
